$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.5747915368188501
$ws.Range("C2").Value = 0.1557691994893844
$ws.Range("E2").Value = 0.1348011902172317
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002379262827167503
$ws.Range("I2").Value = 0.3104406852202857
$ws.Range("M2").Value = 0.2952196453983476
$ws.Range("N2").Value = 0.9219684851426919
$ws.Range("O2").Value = 1.402236214712559
# Row 3
$ws.Range("B3").Value = 0.5025293680304515
$ws.Range("C3").Value = 0.1380220443971041
$ws.Range("E3").Value = 0.1279374319561342
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002381796418243794
$ws.Range("I3").Value = 0.314510296755671
$ws.Range("M3").Value = 0.2625844281756642
$ws.Range("N3").Value = 0.9283761465473219
$ws.Range("O3").Value = 1.398801450834966
# Row 4
$ws.Range("B4").Value = 0.4580514592542784
$ws.Range("C4").Value = 0.1270636041937223
$ws.Range("E4").Value = 0.1238402016742057
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002383434912642701
$ws.Range("I4").Value = 0.317267943082971
$ws.Range("M4").Value = 0.2425994294672904
$ws.Range("N4").Value = 0.9327213886221912
$ws.Range("O4").Value = 1.39792757229398
# Row 5
$ws.Range("B5").Value = 0.4399000338193559
$ws.Range("C5").Value = 0.1225826280945626
$ws.Range("E5").Value = 0.1221997627124409
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002384123510727546
$ws.Range("I5").Value = 0.3184566470400121
$ws.Range("M5").Value = 0.2344687809341224
$ws.Range("N5").Value = 0.9345955401730919
$ws.Range("O5").Value = 1.397881162317475
# Row 6
$ws.Range("B6").Value = 0.436884443608534
$ws.Range("C6").Value = 0.1218376450278527
$ws.Range("E6").Value = 0.1219291265742157
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002384239116165132
$ws.Range("I6").Value = 0.3186579480289282
$ws.Range("M6").Value = 0.2331195033815803
$ws.Range("N6").Value = 0.9349129916484387
$ws.Range("O6").Value = 1.397892135838148
# Row 7
$ws.Range("B7").Value = 0.4578067680051561
$ws.Range("C7").Value = 0.1270032339964189
$ws.Range("E7").Value = 0.1238179601663632
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002383444114652098
$ws.Range("I7").Value = 0.3172837116292051
$ws.Range("M7").Value = 0.2424897223915963
$ws.Range("N7").Value = 0.9327462451631447
$ws.Range("O7").Value = 1.397925693573512
# Row 8
$ws.Range("B8").Value = 0.5498986695469341
$ws.Range("C8").Value = 0.1496628658748023
$ws.Range("E8").Value = 0.132410105242073
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.002380119249950352
$ws.Range("I8").Value = 0.3117900595470005
$ws.Range("M8").Value = 0.2839559998589891
$ws.Range("N8").Value = 0.9240926220670787
$ws.Range("O8").Value = 1.400795038953135
# Row 9
$ws.Range("B9").Value = 0.7295931224206242
$ws.Range("C9").Value = 0.1936039909799945
$ws.Range("E9").Value = 0.1502006397880606
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002374253775980743
$ws.Range("I9").Value = 0.3030783571496123
$ws.Range("M9").Value = 0.3656957737927229
$ws.Range("N9").Value = 0.9103792076494202
$ws.Range("O9").Value = 1.416263544586116
# Row 10
$ws.Range("B10").Value = 0.8610326110449478
$ws.Range("C10").Value = 0.2255818705421007
$ws.Range("E10").Value = 0.1638630695232877
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002370339426194409
$ws.Range("I10").Value = 0.2979444378920135
$ws.Range("M10").Value = 0.4260187915901525
$ws.Range("N10").Value = 0.9022842394805579
$ws.Range("O10").Value = 1.433690014274276
# Row 11
$ws.Range("B11").Value = 0.9206948630979355
$ws.Range("C11").Value = 0.2400623442769643
$ws.Range("E11").Value = 0.1702108771784765
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002368643602384003
$ws.Range("I11").Value = 0.2958860475768503
$ws.Range("M11").Value = 0.4535221902705757
$ws.Range("N11").Value = 0.8990307276335301
$ws.Range("O11").Value = 1.442947587211279
# Row 12
$ws.Range("B12").Value = 0.943267837346184
$ws.Range("C12").Value = 0.2455360409814205
$ws.Range("E12").Value = 0.1726339960954277
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002368013572055072
$ws.Range("I12").Value = 0.295146605249073
$ws.Range("M12").Value = 0.4639460153416053
$ws.Range("N12").Value = 0.8978603192167256
$ws.Range("O12").Value = 1.446645487367761
# Row 13
$ws.Range("B13").Value = 0.9384072426744297
$ws.Range("C13").Value = 0.2443576189065197
$ws.Range("E13").Value = 0.1721112697005225
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002368148721257456
$ws.Range("I13").Value = 0.2953040742960802
$ws.Range("M13").Value = 0.461700663785038
$ws.Range("N13").Value = 0.89810964758626
$ws.Range("O13").Value = 1.4458405120902
# Row 14
$ws.Range("B14").Value = 0.9225523588943929
$ws.Range("C14").Value = 0.2405128649930361
$ws.Range("E14").Value = 0.1704098394986033
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.00236859152639059
$ws.Range("I14").Value = 0.2958244099824263
$ws.Range("M14").Value = 0.454379586867347
$ws.Range("N14").Value = 0.8989332026447059
$ws.Range("O14").Value = 1.443247956937796
# Row 15
$ws.Range("B15").Value = 0.9128381703266086
$ws.Range("C15").Value = 0.2381565674873229
$ws.Range("E15").Value = 0.1693701909834431
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002368864336671359
$ws.Range("I15").Value = 0.2961483485987877
$ws.Range("M15").Value = 0.4498963713776618
$ws.Range("N15").Value = 0.8994456778967148
$ws.Range("O15").Value = 1.441685008780922
# Row 16
$ws.Range("B16").Value = 0.857130827984804
$ws.Range("C16").Value = 0.2246341844130768
$ws.Range("E16").Value = 0.1634509178241998
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002370451954363158
$ws.Range("I16").Value = 0.2980845494351705
$ws.Range("M16").Value = 0.4242226259865589
$ws.Range("N16").Value = 0.9025054906004755
$ws.Range("O16").Value = 1.433111860568914
# Row 17
$ws.Range("B17").Value = 0.8229220397660697
$ws.Range("C17").Value = 0.2163214902224695
$ws.Range("E17").Value = 0.1598538077879397
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.00237144759329911
$ws.Range("I17").Value = 0.2993434382588553
$ws.Range("M17").Value = 0.4084884963472035
$ws.Range("N17").Value = 0.9044924061233033
$ws.Range("O17").Value = 1.428193901228042
# Row 18
$ws.Range("B18").Value = 0.8032338571131277
$ws.Range("C18").Value = 0.2115340077125722
$ws.Range("E18").Value = 0.1577973200079654
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002372028246679916
$ws.Range("I18").Value = 0.3000935845047259
$ws.Range("M18").Value = 0.3994444946180806
$ws.Range("N18").Value = 0.9056756031831625
$ws.Range("O18").Value = 1.42549036588332
# Row 19
$ws.Range("B19").Value = 0.7965657162269508
$ws.Range("C19").Value = 0.209911982438399
$ws.Range("E19").Value = 0.1571031632259476
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002372226219688689
$ws.Range("I19").Value = 0.3003520424698678
$ws.Range("M19").Value = 0.3963833579943383
$ws.Range("N19").Value = 0.9060831496754531
$ws.Range("O19").Value = 1.424596459151303
# Row 20
$ws.Range("B20").Value = 0.8265648935144441
$ws.Range("C20").Value = 0.2172070379933757
$ws.Range("E20").Value = 0.1602354333203451
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.00237134077942103
$ws.Range("I20").Value = 0.2992067282470288
$ws.Range("M20").Value = 0.4101628158814208
$ws.Range("N20").Value = 0.9042767172312551
$ws.Range("O20").Value = 1.428704467277242
# Row 21
$ws.Range("B21").Value = 0.9272098693938347
$ws.Range("C21").Value = 0.2416424284091363
$ws.Range("E21").Value = 0.1709090637841939
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002368461134787146
$ws.Range("I21").Value = 0.2956704870246867
$ws.Range("M21").Value = 0.4565297239092558
$ws.Range("N21").Value = 0.8986896324558273
$ws.Range("O21").Value = 1.444004227600487
# Row 22
$ws.Range("B22").Value = 0.9928711120047637
$ws.Range("C22").Value = 0.257555481689451
$ws.Range("E22").Value = 0.1779977588497133
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002366649863659209
$ws.Range("I22").Value = 0.2935927335432886
$ws.Range("M22").Value = 0.4868849425121908
$ws.Range("N22").Value = 0.8953973309138661
$ws.Range("O22").Value = 1.455124586808608
# Row 23
$ws.Range("B23").Value = 0.9578374533992928
$ws.Range("C23").Value = 0.2490676543130803
$ws.Range("E23").Value = 0.1742039761280836
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002367610119011148
$ws.Range("I23").Value = 0.2946802538196565
$ws.Range("M23").Value = 0.4706790719977221
$ws.Range("N23").Value = 0.8971216456246083
$ws.Range("O23").Value = 1.449086539595669
# Row 24
$ws.Range("B24").Value = 0.8249180249595156
$ws.Range("C24").Value = 0.2168067080106084
$ws.Range("E24").Value = 0.1600628644970143
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002371389044291607
$ws.Range("I24").Value = 0.2992684526397404
$ws.Range("M24").Value = 0.4094058505412335
$ws.Range("N24").Value = 0.9043741027934757
$ws.Range("O24").Value = 1.428473254604825
# Row 25
$ws.Range("B25").Value = 0.6810807146359821
$ws.Range("C25").Value = 0.1817701226800352
$ws.Range("E25").Value = 0.1452852070205282
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.00237577088351317
$ws.Range("I25").Value = 0.3052133976677354
$ws.Range("M25").Value = 0.3435365192057489
$ws.Range("N25").Value = 0.9137409596929302
$ws.Range("O25").Value = 1.411018869802177
